$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto data.
# D-column values look numeric (contain dots) so Excel would otherwise
# auto-convert them to floating point numbers and lose the exact text
# (e.g. trailing zeros, multi-dot "thousands" groupings). Force the cell
# to Text format while writing, then restore the original "Normal" style
# so no stray number formatting is left behind on the cell.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "28.207.83"
$ws.Range("E2").Value = "  +2.15%  "
Set-TextValue $ws.Range("D3") "1.877.90"
$ws.Range("E3").Value = "  +1.63%  "
Set-TextValue $ws.Range("D4") "1.004"
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue $ws.Range("D5") "316.49"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("E6").Value = "  +0.16%  "
Set-TextValue $ws.Range("D7") "0.4318"
$ws.Range("E7").Value = "  +1.33%  "
Set-TextValue $ws.Range("D8") "0.3701"
$ws.Range("E8").Value = "  +1.36%  "
Set-TextValue $ws.Range("D10") "0.8864"
$ws.Range("E10").Value = "  +0.22%  "
Set-TextValue $ws.Range("D11") "21.20"
$ws.Range("E11").Value = "  +2.04%  "
Set-TextValue $ws.Range("D12") "1.895.24"
$ws.Range("E12").Value = "  -1.17%  "
Set-TextValue $ws.Range("D13") "5.498"
$ws.Range("E13").Value = "  +2.72%  "
Set-TextValue $ws.Range("D14") "6.624"
$ws.Range("E14").Value = "  +1.19%  "
Set-TextValue $ws.Range("D15") "0.06973"
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("E16").Value = "  +0.20%  "
Set-TextValue $ws.Range("D17") "81.23"
$ws.Range("E17").Value = "  +2.86%  "
Set-TextValue $ws.Range("D18") "0.000009143"
$ws.Range("E18").Value = "  +2.80%  "
Set-TextValue $ws.Range("D19") "1.004"
$ws.Range("E19").Value = "  +0.21%  "
Set-TextValue $ws.Range("D20") "15.63"
$ws.Range("E20").Value = "  +1.18%  "
Set-TextValue $ws.Range("D21") "28.164.11"
$ws.Range("E21").Value = "  +1.95%  "
Set-TextValue $ws.Range("D22") "5.091"
$ws.Range("E22").Value = "  +2.14%  "
Set-TextValue $ws.Range("D23") "10.96"
$ws.Range("E23").Value = "  +2.73%  "
Set-TextValue $ws.Range("D24") "2.133.04"
$ws.Range("E24").Value = "  +1.44%  "
Set-TextValue $ws.Range("D25") "1.980"
$ws.Range("E25").Value = "  +0.46%  "
Set-TextValue $ws.Range("D26") "154.54"
$ws.Range("E26").Value = "  +0.56%  "
Set-TextValue $ws.Range("D27") "18.74"
$ws.Range("E27").Value = "  -1.22%  "
Set-TextValue $ws.Range("D28") "5.435"
$ws.Range("E28").Value = "  +3.31%  "
Set-TextValue $ws.Range("D29") "118.53"
$ws.Range("E29").Value = "  -2.85%  "
Set-TextValue $ws.Range("D30") "1.904"
$ws.Range("E30").Value = "  -0.67%  "
Set-TextValue $ws.Range("D31") "0.08982"
$ws.Range("E31").Value = "  +0.41%  "
Set-TextValue $ws.Range("D32") "0.7952"
$ws.Range("E32").Value = "  +4.19%  "
Set-TextValue $ws.Range("D33") "4.698"
$ws.Range("E33").Value = "  +2.63%  "
Set-TextValue $ws.Range("D34") "1.172"
$ws.Range("E34").Value = "  +6.24%  "
Set-TextValue $ws.Range("D35") "2.979"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +3.65%  "
$ws.Range("E37").Value = "  +0.24%  "
Set-TextValue $ws.Range("D38") "0.05482"
$ws.Range("E38").Value = "  +1.95%  "
Set-TextValue $ws.Range("D39") "0.01969"
$ws.Range("E39").Value = "  +1.00%  "
Set-TextValue $ws.Range("D40") "2.887"
$ws.Range("E40").Value = "  +2.43%  "
Set-TextValue $ws.Range("D41") "0.1698"
$ws.Range("E41").Value = "  +2.56%  "
Set-TextValue $ws.Range("D42") "0.5177"
$ws.Range("E42").Value = "  +1.12%  "
Set-TextValue $ws.Range("D43") "6.888"
$ws.Range("E43").Value = "  -0.72%  "
Set-TextValue $ws.Range("D44") "8.581"
Set-TextValue $ws.Range("D45") "10.60"
$ws.Range("E45").Value = "  +2.05%  "
Set-TextValue $ws.Range("D46") "0.06596"
$ws.Range("E46").Value = "  +0.29%  "
Set-TextValue $ws.Range("D47") "0.4762"
$ws.Range("E47").Value = "  -0.07%  "
Set-TextValue $ws.Range("D48") "105.75"
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("E49").Value = "  +0.23%  "
Set-TextValue $ws.Range("D50") "1.661"
$ws.Range("E50").Value = "  +1.85%  "
Set-TextValue $ws.Range("D51") "1.841"
$ws.Range("E51").Value = "  +4.37%  "
